$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.936.73'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.637.54'
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.29'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("E10").Value = '  +0.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0883'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").Value = '1.869.27'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '1.639.55'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '27.940.86'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '231.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -0.27%  '
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.96%  '
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("E24").Value = '  -4.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.52%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E28").Value = '  -0.73%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.10%  '
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("D33").Value = '1.401.55'
$ws.Range("E33").Value = '  -3.87%  '
$ws.Range("E34").Value = '  -1.47%  '
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("E38").Value = '  -0.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.926'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.874'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '66.93'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = '1.778.34'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1000'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.59'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.52%  '
